# "upd theory of discrete math + rename comp system to CS and discrete math to DM"
#
# The subjects sheet lists one course per row in column A. This adds a short
# abbreviation in column B next to the two affected courses:
#   - "вычислительные системы" (computer/computational systems) -> "CS"
#   - "дискретная математика"  (discrete math)                  -> "DM"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Row 4: вычислительные системы -> tag as "CS"
$ws.Range("B4").Value = "CS"

# Row 8: дискретная математика -> tag as "DM"
$ws.Range("B8").Value = "DM"

# Match the author's resulting view/selection state (scrolled down, B9 active)
$ws.Select()
$ws.Range("B9").Select()
try {
    $excel.ActiveWindow.ScrollRow = 3
} catch {
    # view-scroll state isn't always settable in headless hosts; non-fatal
}
